$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.877
$ws.Range("A8").Value = -22.305
$ws.Range("A10").Value = -21.871
$ws.Range("A12").Value = -21.474
$ws.Range("D12").Value = -6.873
$ws.Range("D15").Value = -8.311999999999999
$ws.Range("D17").Value = -8.193
$ws.Range("A18").Value = -22.095
$ws.Range("D26").Value = -7.255000000000001
$ws.Range("D27").Value = -7.728
$ws.Range("D28").Value = -8.047999999999998
$ws.Range("A37").Value = -20.458
$ws.Range("D37").Value = -8.257000000000001
$ws.Range("D47").Value = -7.452
$ws.Range("A55").Value = -22.204
$ws.Range("D65").Value = -7.741
$ws.Range("A68").Value = -21.53
$ws.Range("D73").Value = -8.129999999999999
$ws.Range("A77").Value = -20.522
$ws.Range("A78").Value = -19.951
$ws.Range("A81").Value = -21.72
$ws.Range("A82").Value = -22.209
$ws.Range("D84").Value = -7.994
$ws.Range("D85").Value = -8.559999999999999
$ws.Range("D93").Value = -7.007000000000001
$ws.Range("D95").Value = -7.557
$ws.Range("D98").Value = -7.203999999999999
$ws.Range("D99").Value = -8.189
$ws.Range("D101").Value = -8.214
